$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("by Coach")

# Toggle "Started" (column C) Yes/No values for the rows that changed.
$rowsToToggle = 3,9,31,32,42,45,51,54,55,56,76,77
foreach ($r in $rowsToToggle) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value()
    if ($current -eq "Yes") {
        $cell.Value = "No"
    } else {
        $cell.Value = "Yes"
    }
}

# Update the view: the active selection moved down towards the bottom
# of the data as the user scrolled/edited (sheet stays frozen on row 1).
$ws.Activate()
$ws.Range("C84").Select()
